$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename: Nrc_Id -> NRC_ID
$ws.Cells.Item(1, 1).Value = "NRC_ID"

# 2. Widen column C (Ref_Task) to match the report's new layout
$ws.Columns.Item(3).ColumnWidth = 24.71

# 3. Append a new data row (row 8) by duplicating row 7 (same B-HNU / SWC-28778-01-01 /
#    ZONE C ... complaint) then updating the NRC id and total man-hours.
$ws.Range("A7:K7").Copy()
$ws.Range("A8:K8").PasteSpecial(-4163)
$ws.Cells.Item(8, 1).Value = "QZL0070"
$ws.Cells.Item(8, 10).Value = 50

# 4. Bump the sheet's outline-level bookkeeping from 6 to 7 (group/ungroup a
#    throwaway row beyond the data range so no stray outline attributes remain).
for ($i = 0; $i -lt 7; $i++) {
    $ws.Rows.Item(9).Group()
}
$ws.Rows.Item(9).Delete()

# 5. Leave the selection on A2, matching the saved view state.
[void]$ws.Range("A2").Select()
